$wb = $excel.ActiveWorkbook

# Rename "rooms" sheet to "beds"
$wsBeds = $wb.Worksheets.Item("rooms")
$wsBeds.Name = "beds"

# Update the header row terminology from "rooms" to "beds"
$wsBeds.Range("A1").Value = "all_beds"
$wsBeds.Range("B1").Value = "new_beds"
$wsBeds.Range("C1").Value = "old_beds"
$wsBeds.Range("E1").Value = "new_beds_service"
$wsBeds.Range("F1").Value = "old_beds_service"
$wsBeds.Range("G1").Value = "beds_capacities"

# Make the "beds" sheet the active tab, and change its selection
$wsBeds.Activate()
$wsBeds.Range("C22").Select()
